$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "18:06 03-Dec-23"
$ws.Range("C11").Value = "Ẩn danh"
$ws.Range("D11").Value = "ccccccccccc"
